$d = $word.ActiveDocument

$d.Content.Find.Execute("691÷7=98, 5", $true, $false, $false, $false, $false, $true, 1, $false, "104÷2=52, 0", 2) | Out-Null
$d.Content.Find.Execute("825÷5=165, 0", $true, $false, $false, $false, $false, $true, 1, $false, "294÷4=73, 2", 2) | Out-Null
$d.Content.Find.Execute("964÷9=107, 1", $true, $false, $false, $false, $false, $true, 1, $false, "132÷7=18, 6", 2) | Out-Null
$d.Content.Find.Execute("211÷7=30, 1", $true, $false, $false, $false, $false, $true, 1, $false, "802÷6=133, 4", 2) | Out-Null
$d.Content.Find.Execute("971÷5=194, 1", $true, $false, $false, $false, $false, $true, 1, $false, "382÷5=76, 2", 2) | Out-Null
$d.Content.Find.Execute("515÷8=64, 3", $true, $false, $false, $false, $false, $true, 1, $false, "643÷7=91, 6", 2) | Out-Null
$d.Content.Find.Execute("671÷8=83, 7", $true, $false, $false, $false, $false, $true, 1, $false, "598÷4=149, 2", 2) | Out-Null
$d.Content.Find.Execute("918÷3=306, 0", $true, $false, $false, $false, $false, $true, 1, $false, "297÷5=59, 2", 2) | Out-Null
$d.Content.Find.Execute("165÷7=23, 4", $true, $false, $false, $false, $false, $true, 1, $false, "399÷3=133, 0", 2) | Out-Null
$d.Content.Find.Execute("920÷9=102, 2", $true, $false, $false, $false, $false, $true, 1, $false, "645÷7=92, 1", 2) | Out-Null
$d.Content.Find.Execute("522÷6=87, 0", $true, $false, $false, $false, $false, $true, 1, $false, "507÷7=72, 3", 2) | Out-Null
$d.Content.Find.Execute("140÷8=17, 4", $true, $false, $false, $false, $false, $true, 1, $false, "718÷7=102, 4", 2) | Out-Null
$d.Content.Find.Execute("409÷7=58, 3", $true, $false, $false, $false, $false, $true, 1, $false, "587÷3=195, 2", 2) | Out-Null
$d.Content.Find.Execute("686÷4=171, 2", $true, $false, $false, $false, $false, $true, 1, $false, "991÷8=123, 7", 2) | Out-Null
$d.Content.Find.Execute("518÷2=259, 0", $true, $false, $false, $false, $false, $true, 1, $false, "281÷8=35, 1", 2) | Out-Null
$d.Content.Find.Execute("390÷4=97, 2", $true, $false, $false, $false, $false, $true, 1, $false, "469÷2=234, 1", 2) | Out-Null
$d.Content.Find.Execute("939÷3=313, 0", $true, $false, $false, $false, $false, $true, 1, $false, "627÷3=209, 0", 2) | Out-Null
$d.Content.Find.Execute("515÷6=85, 5", $true, $false, $false, $false, $false, $true, 1, $false, "885÷6=147, 3", 2) | Out-Null
$d.Content.Find.Execute("311÷2=155, 1", $true, $false, $false, $false, $false, $true, 1, $false, "442÷2=221, 0", 2) | Out-Null
$d.Content.Find.Execute("382÷4=95, 2", $true, $false, $false, $false, $false, $true, 1, $false, "564÷6=94, 0", 2) | Out-Null
$d.Content.Find.Execute("866÷4=216, 2", $true, $false, $false, $false, $false, $true, 1, $false, "164÷5=32, 4", 2) | Out-Null
$d.Content.Find.Execute("966÷9=107, 3", $true, $false, $false, $false, $false, $true, 1, $false, "883÷6=147, 1", 2) | Out-Null
$d.Content.Find.Execute("585÷3=195, 0", $true, $false, $false, $false, $false, $true, 1, $false, "748÷4=187, 0", 2) | Out-Null
$d.Content.Find.Execute("886÷7=126, 4", $true, $false, $false, $false, $false, $true, 1, $false, "633÷4=158, 1", 2) | Out-Null
$d.Content.Find.Execute("298÷4=74, 2", $true, $false, $false, $false, $false, $true, 1, $false, "899÷7=128, 3", 2) | Out-Null
